$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Diameter_mm" header in N1, mirroring the format of M1 ("Diameter")
$ws.Range("N1").Value = "Diameter_mm"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122) | Out-Null

# Fill in the metric (mm) conversion formulas for N2:N59, converting the
# inches diameter already computed in column M (M2:M59) to millimetres -
# matching the way the existing M column itself is split into a single
# formula on the first row and a shared formula for the remaining rows.
$ws.Range("N2").Formula = "=M2*25.4"
$ws.Range("N3:N59").Formula = "=M3*25.4"
